# CSV file upload was integrated: update the registration number on the
# student psychomotor sheet and move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Registration number changed (new CSV upload content)
$ws.Range("A2").Value = "2021C123450003"

# Active cell/selection moved from H2 to A2
$ws.Range("A2").Select()
